$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("G5").Value = 2.15
$ws.Range("I5").Value = 3.25
$ws.Range("Z5").Value = 21
$ws.Range("AO5").Value = 12
$ws.Range("AX5").Value = 17
$ws.Range("AY5").Value = 23
$ws.Range("BA5").Value = 67

# Row 6 updates
$ws.Range("N6").Value = 10

# Row 9 updates
$ws.Range("Q9").Value = 2.5
$ws.Range("R9").Value = 1.53
